$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change R1's text from "Identificacion" to "Número ID" (same as E1)
$ws.Range("R1").Value = "Número ID"

# Add a new row 2 (A2:R2) - empty cells with a thin border around each cell
$rng = $ws.Range("A2:R2")
$rng.Borders.Color = 0.4
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2
